$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $b, $c, $d, $dIsNumericLooking, $e) {
    if ($b -ne $null) { $ws.Cells.Item($row, 2).Value = $b }
    if ($c -ne $null) { $ws.Cells.Item($row, 3).Value = $c }
    if ($d -ne $null) {
        if ($dIsNumericLooking -eq $true) {
            $ws.Cells.Item($row, 4).NumberFormat = "@"
        }
        $ws.Cells.Item($row, 4).Value = $d
    }
    if ($e -ne $null) { $ws.Cells.Item($row, 5).Value = $e }
}

# Row 2 - Bitcoin
Set-Row 2 $null $null "55.147.32" $false "  -1.58%  "
# Row 3 - Ethereum
Set-Row 3 $null $null "2.344.49" $false "  -4.82%  "
# Row 4 - TetherUSD
Set-Row 4 $null $null $null $false "  -0.05%  "
# Row 5 - BNB
Set-Row 5 $null $null "475.58" $true "  -2.41%  "
# Row 6 - Solana
Set-Row 6 $null $null "145.05" $true "  -0.03%  "
# Row 7 - USDC
Set-Row 7 $null $null "0.999" $true "  +0.04%  "
# Row 8 - XRP
Set-Row 8 $null $null "0.611" $true "  +20.31%  "
# Row 9 - LidoStakedEther
Set-Row 9 $null $null "2.343.86" $false "  -4.97%  "
# Row 10 - Dogecoin
Set-Row 10 $null $null $null $false "  -1.12%  "
# Row 11 - Toncoin
Set-Row 11 $null $null "5.44" $true "  -6.45%  "
# Row 12 - Cardano
Set-Row 12 $null $null "0.325" $true "  -1.71%  "
# Row 13 - TRON
Set-Row 13 $null $null $null $false "  +1.33%  "
# Row 14 - WrappedliquidstakedEther2.0
Set-Row 14 $null $null "2.747.69" $false "  -4.98%  "
# Row 15 - WrappedBTC
Set-Row 15 $null $null "55.139.91" $false "  -1.63%  "
# Row 16 - Avalanche
Set-Row 16 $null $null "19.93" $true "  -5.39%  "
# Row 17 - ShibaInu
Set-Row 17 $null $null $null $false "  -4.91%  "
# Row 18 - WrappedEther
Set-Row 18 $null $null "2.348.33" $false "  -5.01%  "
# Row 19 - Polkadot
Set-Row 19 $null $null $null $false "  +0.94%  "
# Row 20 - BitcoinCash
Set-Row 20 $null $null "313.73" $true "  -0.90%  "
# Row 21 - Chainlink
Set-Row 21 $null $null $null $false "  -4.88%  "
# Row 22 - Dai
Set-Row 22 $null $null $null $false "  +0.20%  "
# Row 23 - Uniswap
Set-Row 23 $null $null "5.64" $true "  -2.44%  "
# Row 24 - Litecoin
Set-Row 24 $null $null "56.14" $true "  -3.93%  "
# Row 25 - Binance-PegBSC-USD
Set-Row 25 $null $null "0.999" $true "  -0.04%  "
# Row 26 - Polygon
Set-Row 26 $null $null $null $false "  -4.19%  "
# Row 27 - Kaspa
Set-Row 27 $null $null $null $false "  -5.42%  "
# Row 28 - WrappedeETH
Set-Row 28 $null $null "2.437.97" $false "  -5.35%  "
# Row 29 - InternetComputer(DFINITY)
Set-Row 29 $null $null "7.05" $true "  -8.26%  "
# Row 30 - USDe (unchanged)
# Row 31 - PEPE
Set-Row 31 $null $null "0.0₃0740" $false "  -5.45%  "
# Row 32 - Monero
Set-Row 32 $null $null "145.76" $true "  -1.35%  "
# Row 33 - EthereumClassic
Set-Row 33 $null $null "18.07" $true "  -0.74%  "
# Row 34 - PancakeSwap
Set-Row 34 $null $null $null $false "  -2.11%  "
# Row 35 - Aptos
Set-Row 35 $null $null "5.08" $true "  -1.67%  "

# Rows 36/37 swap: NEARProtocol <-> ImmutableX
Set-Row 36 "ImmutableX" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" "1.09" $true "  -3.97%  "
Set-Row 37 "NEARProtocol" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near" "3.59" $true "  -3.44%  "

# Row 38 - Fetch.AI
Set-Row 38 $null $null $null $false "  -5.98%  "
# Row 39 - Stellar
Set-Row 39 $null $null $null $false "  +9.88%  "
# Row 40 - OKB
Set-Row 40 $null $null "33.58" $true "  -0.94%  "
# Row 41 - FirstDigitalUSD
Set-Row 41 $null $null "0.997" $true "  +0.17%  "
# Row 42 - Stacks
Set-Row 42 $null $null $null $false "  -0.22%  "
# Row 43 - Filecoin
Set-Row 43 $null $null $null $false "  -4.09%  "
# Row 44 - Mantle
Set-Row 44 $null $null $null $false "  -4.38%  "

# Rows 45/46 swap: Hedera <-> WhiteBITCoin
Set-Row 45 "WhiteBITCoin" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt" "10.15" $true "  -0.48%  "
Set-Row 46 "Hedera" "https://coinranking.com/coin/jad286TjB+hedera-hbar" "0.0513" $true "  -7.09%  "

# Row 47 - Bittensor
Set-Row 47 $null $null "248.86" $true "  -4.72%  "
# Row 48 - VeChain
Set-Row 48 $null $null $null $false "  -2.99%  "
# Row 49 - RenderToken
Set-Row 49 $null $null $null $false "  -7.79%  "
# Row 50 - Maker
Set-Row 50 $null $null "1.788.70" $false "  -4.40%  "
# Row 51 - EnergySwap
Set-Row 51 $null $null $null $false "  -5.43%  "
